$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S11").Value = -0.0154

$ws.Range("G22").Value = -0.9172
$ws.Range("H22").Value = -0.2334
$ws.Range("I22").Value = -0.8409
$ws.Range("J22").Value = -0.4083

$ws.Range("K23").Value = 0.0142
$ws.Range("L23").Value = 0.124
$ws.Range("M23").Value = 0.1722
$ws.Range("N23").Value = -0.1072
$ws.Range("O23").Value = -0.1102
$ws.Range("P23").Value = -0.095
$ws.Range("Q23").Value = -0.08
$ws.Range("R23").Value = -0.0445
$ws.Range("S23").Value = -0.0796

$ws.Range("G30").Value = -3.8339
$ws.Range("H30").Value = -4.8217
$ws.Range("I30").Value = -2.4774
$ws.Range("J30").Value = -1.0032

$ws.Range("K31").Value = -1.5351
$ws.Range("L31").Value = -1.6397
$ws.Range("M31").Value = -0.574
$ws.Range("N31").Value = -0.5699
$ws.Range("O31").Value = -0.6183
$ws.Range("P31").Value = -0.1956
$ws.Range("Q31").Value = -0.1398
$ws.Range("R31").Value = 0.0254
$ws.Range("S31").Value = -0.052

$ws.Range("S67").Value = 0.0312

$ws.Range("G78").Value = 0.0092
$ws.Range("H78").Value = 0.0089
$ws.Range("I78").Value = 0.0086
$ws.Range("J78").Value = 0.0205

$ws.Range("K79").Value = -0.0038
$ws.Range("L79").Value = -0.0161
$ws.Range("M79").Value = -0.0274
$ws.Range("N79").Value = -0.0558
$ws.Range("O79").Value = -0.0439
$ws.Range("P79").Value = -0.0312
$ws.Range("Q79").Value = -0.0184
$ws.Range("R79").Value = -0.0036
$ws.Range("S79").Value = -0.0013

$ws.Range("G86").Value = 0.0092
$ws.Range("H86").Value = 0.0089
$ws.Range("I86").Value = 0.0086
$ws.Range("J86").Value = 0.0205

$ws.Range("K87").Value = -0.0632
$ws.Range("L87").Value = 0.1503
$ws.Range("M87").Value = 0.3096
$ws.Range("N87").Value = 0.1331
$ws.Range("O87").Value = 0.1242
$ws.Range("P87").Value = 0.1485
$ws.Range("Q87").Value = 0.196
$ws.Range("R87").Value = 0.1191
$ws.Range("S87").Value = -0.0862
